$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.063063333333333
$ws.Range("H2").Value = 6.189190000000001
$ws.Range("I2").Value = 0.1875199417503197
$ws.Range("J2").Value = 0.1875199417503197
$ws.Range("M2").Value = 14.25737566666667
$ws.Range("N2").Value = 42.772127
$ws.Range("O2").Value = 0.2087950866344732
$ws.Range("P2").Value = 0.2087950866344732
$ws.Range("Q2").Value = 29.41386896745889
$ws.Range("R2").Value = 264.72482070713
$ws.Range("S2").Value = 0.03915324248344936
$ws.Range("T2").Value = 0.03915324248344937

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.063063333333333
$ws.Range("H3").Value = 6.189190000000001
$ws.Range("I3").Value = 0.1875199417503197
$ws.Range("J3").Value = 0.1875199417503197
$ws.Range("N3").Value = 87.128332
$ws.Range("O3").Value = 0.4253229592313036
$ws.Range("P3").Value = 0.4253229592313036
$ws.Range("Q3").Value = 59.91708901456445
$ws.Range("R3").Value = 539.2538011310801
$ws.Range("S3").Value = 0.07975653654012764
$ws.Range("T3").Value = 0.07975653654012765

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.063063333333333
$ws.Range("H4").Value = 6.189190000000001
$ws.Range("I4").Value = 0.1875199417503197
$ws.Range("J4").Value = 0.1875199417503197
$ws.Range("M4").Value = 20.11084633333333
$ws.Range("N4").Value = 60.332539
$ws.Range("O4").Value = 0.2945174484164121
$ws.Range("P4").Value = 0.2945174484164122
$ws.Range("Q4").Value = 41.48994967260111
$ws.Range("R4").Value = 373.4095470534101
$ws.Range("S4").Value = 0.05522789477149838
$ws.Range("T4").Value = 0.0552278947714984

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.063063333333333
$ws.Range("H5").Value = 6.189190000000001
$ws.Range("I5").Value = 0.1875199417503197
$ws.Range("J5").Value = 0.1875199417503197
$ws.Range("M5").Value = 4.873057999999999
$ws.Range("N5").Value = 14.619174
$ws.Range("O5").Value = 0.07136450571781097
$ws.Range("P5").Value = 0.07136450571781099
$ws.Range("Q5").Value = 10.05342728100667
$ws.Range("R5").Value = 90.48084552906
$ws.Range("S5").Value = 0.01338226795524427
$ws.Range("T5").Value = 0.01338226795524427

# Row 6
$ws.Range("I6").Value = 0.5238509469163369
$ws.Range("J6").Value = 0.5238509469163369
$ws.Range("M6").Value = 14.25737566666667
$ws.Range("N6").Value = 42.772127
$ws.Range("O6").Value = 0.2087950866344732
$ws.Range("P6").Value = 0.2087950866344732
$ws.Range("Q6").Value = 82.16983733704755
$ws.Range("R6").Value = 739.5285360334279
$ws.Range("S6").Value = 0.1093775038449474
$ws.Range("T6").Value = 0.1093775038449474

# Row 7
$ws.Range("I7").Value = 0.5238509469163369
$ws.Range("J7").Value = 0.5238509469163369
$ws.Range("N7").Value = 87.128332
$ws.Range("O7").Value = 0.4253229592313036
$ws.Range("P7").Value = 0.4253229592313036
$ws.Range("Q7").Value = 167.3828581844498
$ws.Range("S7").Value = 0.2228058349385769
$ws.Range("T7").Value = 0.222805834938577

# Row 8
$ws.Range("I8").Value = 0.5238509469163369
$ws.Range("J8").Value = 0.5238509469163369
$ws.Range("M8").Value = 20.11084633333333
$ws.Range("N8").Value = 60.332539
$ws.Range("O8").Value = 0.2945174484164121
$ws.Range("P8").Value = 0.2945174484164122
$ws.Range("Q8").Value = 115.9052697042884
$ws.Range("R8").Value = 1043.147427338596
$ws.Range("S8").Value = 0.1542832442363209
$ws.Range("T8").Value = 0.1542832442363209

# Row 9
$ws.Range("I9").Value = 0.5238509469163369
$ws.Range("J9").Value = 0.5238509469163369
$ws.Range("M9").Value = 4.873057999999999
$ws.Range("N9").Value = 14.619174
$ws.Range("O9").Value = 0.07136450571781097
$ws.Range("P9").Value = 0.07136450571781099
$ws.Range("Q9").Value = 28.08499912997066
$ws.Range("R9").Value = 252.764992169736
$ws.Range("S9").Value = 0.03738436389649162
$ws.Range("T9").Value = 0.03738436389649163

# Row 10
$ws.Range("G10").Value = 2.101774
$ws.Range("H10").Value = 6.305322
$ws.Range("I10").Value = 0.1910385065181404
$ws.Range("J10").Value = 0.1910385065181404
$ws.Range("M10").Value = 14.25737566666667
$ws.Range("N10").Value = 42.772127
$ws.Range("O10").Value = 0.2087950866344732
$ws.Range("P10").Value = 0.2087950866344732
$ws.Range("Q10").Value = 29.96578148443268
$ws.Range("R10").Value = 269.692033359894
$ws.Range("S10").Value = 0.03988790151897548
$ws.Range("T10").Value = 0.03988790151897549

# Row 11
$ws.Range("G11").Value = 2.101774
$ws.Range("H11").Value = 6.305322
$ws.Range("I11").Value = 0.1910385065181404
$ws.Range("J11").Value = 0.1910385065181404
$ws.Range("N11").Value = 87.128332
$ws.Range("O11").Value = 0.4253229592313036
$ws.Range("P11").Value = 0.4253229592313036
$ws.Range("Q11").Value = 61.04135428698934
$ws.Range("R11").Value = 549.3721885829041
$ws.Range("S11").Value = 0.08125306291942413
$ws.Range("T11").Value = 0.08125306291942415

# Row 12
$ws.Range("G12").Value = 2.101774
$ws.Range("H12").Value = 6.305322
$ws.Range("I12").Value = 0.1910385065181404
$ws.Range("J12").Value = 0.1910385065181404
$ws.Range("M12").Value = 20.11084633333333
$ws.Range("N12").Value = 60.332539
$ws.Range("O12").Value = 0.2945174484164121
$ws.Range("P12").Value = 0.2945174484164122
$ws.Range("Q12").Value = 42.26845394139534
$ws.Range("R12").Value = 380.416085472558
$ws.Range("S12").Value = 0.05626417348900482
$ws.Range("T12").Value = 0.05626417348900483

# Row 13
$ws.Range("G13").Value = 2.101774
$ws.Range("H13").Value = 6.305322
$ws.Range("I13").Value = 0.1910385065181404
$ws.Range("J13").Value = 0.1910385065181404
$ws.Range("M13").Value = 4.873057999999999
$ws.Range("N13").Value = 14.619174
$ws.Range("O13").Value = 0.07136450571781097
$ws.Range("P13").Value = 0.07136450571781099
$ws.Range("Q13").Value = 10.242066604892
$ws.Range("R13").Value = 92.178599444028
$ws.Range("S13").Value = 0.0136333685907359
$ws.Range("T13").Value = 0.0136333685907359

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.073675666666667
$ws.Range("H14").Value = 3.221027
$ws.Range("I14").Value = 0.0975906048152031
$ws.Range("J14").Value = 0.09759060481520311
$ws.Range("M14").Value = 14.25737566666667
$ws.Range("N14").Value = 42.772127
$ws.Range("O14").Value = 0.2087950866344732
$ws.Range("P14").Value = 0.2087950866344732
$ws.Range("Q14").Value = 15.30779732382545
$ws.Range("R14").Value = 137.770175914429
$ws.Range("S14").Value = 0.02037643878710097
$ws.Range("T14").Value = 0.02037643878710097

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.073675666666667
$ws.Range("H15").Value = 3.221027
$ws.Range("I15").Value = 0.0975906048152031
$ws.Range("J15").Value = 0.09759060481520311
$ws.Range("N15").Value = 87.128332
$ws.Range("O15").Value = 0.4253229592313036
$ws.Range("P15").Value = 0.4253229592313036
$ws.Range("Q15").Value = 31.18252331521822
$ws.Range("R15").Value = 280.642709836964
$ws.Range("S15").Value = 0.04150752483317489
$ws.Range("T15").Value = 0.0415075248331749

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.073675666666667
$ws.Range("H16").Value = 3.221027
$ws.Range("I16").Value = 0.0975906048152031
$ws.Range("J16").Value = 0.09759060481520311
$ws.Range("M16").Value = 20.11084633333333
$ws.Range("N16").Value = 60.332539
$ws.Range("O16").Value = 0.2945174484164121
$ws.Range("P16").Value = 0.2945174484164122
$ws.Range("Q16").Value = 21.59252634417255
$ws.Range("R16").Value = 194.332737097553
$ws.Range("S16").Value = 0.02874213591958804
$ws.Range("T16").Value = 0.02874213591958805

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.073675666666667
$ws.Range("H17").Value = 3.221027
$ws.Range("I17").Value = 0.0975906048152031
$ws.Range("J17").Value = 0.09759060481520311
$ws.Range("M17").Value = 4.873057999999999
$ws.Range("N17").Value = 14.619174
$ws.Range("O17").Value = 0.07136450571781097
$ws.Range("P17").Value = 0.07136450571781099
$ws.Range("Q17").Value = 5.232083796855333
$ws.Range("R17").Value = 47.088754171698
$ws.Range("S17").Value = 0.006964505275339193
$ws.Range("T17").Value = 0.006964505275339194

